$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style used by ordinary (unstyled) data cells, so that we can
# restore it after temporarily forcing a cell to Text format (needed so
# Excel does not "helpfully" reinterpret numeric-looking strings as numbers).
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.913.19"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.205.17"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.39"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +3.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.08"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +2.32%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.199.08"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  -0.78%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.06"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.12"
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = "  +4.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.732.03"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.963.33"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.43"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +3.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.203.97"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  +1.34%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.26"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.32"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  +3.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.743"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  +2.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.30"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.04"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.99"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.38"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +3.99%  "
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("E29").Value = "  +3.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.87"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.85"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +8.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.35"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +2.09%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.56"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.11"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0908"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  +0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "484.52"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0420"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -3.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.86"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +2.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.299"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +4.81%  "
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0648"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  +8.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.944.10"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  -3.90%  "
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.38"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  +1.09%  "
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.39"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  -0.59%  "
